# Updated Framework to run cases individually.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
$ws.Activate()

$ws.Range("B4").Value = "Y"
$ws.Range("A4").Select()
